$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processing Initialization")

# Group "Aged 1 to 20" (rows 2-20): conveyor index was counting up 1..19, flip it to count down 19..1
$ws.Range("B2").Value = 19
$ws.Range("B3").Formula = "=B2-1"
$ws.Range("B4:B20").Formula = "=B3-1"

# Group "Aged 20 to 40" (rows 21-40): conveyor index was counting up 1..20, flip it to count down 20..1
$ws.Range("B21").Value = 20
$ws.Range("B22").Formula = "=B21-1"
$ws.Range("B23:B40").Formula = "=B22-1"

# Groups "Aged 40 to 60" / "Aged 60 to 80" / "Aged Over 80" reuse the formulas above via
# B41=B21, B42:B100 shared "=B22" (offset -20), so they recompute automatically.

# Re-select the sheet/cell that was active when the file was last saved.
$ws.Activate()
$ws.Range("B41").Select()
